# Regenerate merged AHB files
# 1) Rename header labels from *_old / *_new to *_FV2410 / *_FV2504
# 2) Wrap the data range in an Excel Table (ListObject)
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffixMap = @{
    "Segmentname_old"        = "Segmentname_FV2410"
    "Segmentgruppe_old"      = "Segmentgruppe_FV2410"
    "Segment_old"            = "Segment_FV2410"
    "Datenelement_old"       = "Datenelement_FV2410"
    "Segment ID_old"         = "Segment ID_FV2410"
    "Code_old"               = "Code_FV2410"
    "Qualifier_old"          = "Qualifier_FV2410"
    "Beschreibung_old"       = "Beschreibung_FV2410"
    "Bedingungsausdruck_old" = "Bedingungsausdruck_FV2410"
    "Bedingung_old"          = "Bedingung_FV2410"
    "Segmentname_new"        = "Segmentname_FV2504"
    "Segmentgruppe_new"      = "Segmentgruppe_FV2504"
    "Segment_new"            = "Segment_FV2504"
    "Datenelement_new"       = "Datenelement_FV2504"
    "Segment ID_new"         = "Segment ID_FV2504"
    "Code_new"               = "Code_FV2504"
    "Qualifier_new"          = "Qualifier_FV2504"
    "Beschreibung_new"       = "Beschreibung_FV2504"
    "Bedingungsausdruck_new" = "Bedingungsausdruck_FV2504"
    "Bedingung_new"          = "Bedingung_FV2504"
}

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value2
    if ($suffixMap.ContainsKey($old)) {
        $cell.Value = $suffixMap[$old]
    }
}

# Preserve the header row's existing look (bold / fill / border / center+wrap):
# stash it, strip it so the table-creation step below doesn't freeze it into a
# dxf override, then paste it straight back once the table exists.
$headerRange = $ws.Range("A1:U1")
$headerRange.Copy()
$headerRange.ClearFormats()

# Wrap the used range in a table
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$headerRange.PasteSpecial(-4122)

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
